# Generate Report for Handoff
# Update the status/handoff info for the "d44fcb6e-9b32-40a8-b504-7e89e1dfe160.md" row
# (row 3) across the Overview, zh-cn and de-de sheets: it has moved from
# "In Translation" to "Ready for handoff", with fresh handoff timestamps.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-03-24 02:19:54"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-24 02:19:51"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-24 02:19:54"
